# Refresh the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped values, mirroring the GitHub Actions "Updated
# cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some Price values (e.g. "550.68", "1.00") look like plain decimal
# numbers to Excel's auto-detection and would otherwise be silently coerced
# into numeric cells (losing the trailing zero / exact text). Force the
# cell to Text format, assign the literal string, then restore the
# "Normal" cell style so no stray number-format/style survives on the
# cell itself (matches the original inline-string text cells).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "61.895.43"
$ws.Range("E2").Value = "  -2.32%  "

$ws.Range("D3").Value = "2.578.80"

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue "D5" "550.68"
$ws.Range("E5").Value = "  -0.43%  "

Set-TextValue "D6" "154.74"
$ws.Range("E6").Value = "  -2.34%  "

$ws.Range("E7").Value = "  +0.00%  "

Set-TextValue "D8" "0.592"
$ws.Range("E8").Value = "  +1.28%  "

$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("E10").Value = "  -0.92%  "

Set-TextValue "D11" "5.49"
$ws.Range("E11").Value = "  +2.88%  "

$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("D13").Value = "3.031.38"
$ws.Range("E13").Value = "  -3.95%  "

Set-TextValue "D14" "25.43"
$ws.Range("E14").Value = "  -3.58%  "

$ws.Range("D15").Value = "61.824.26"
$ws.Range("E15").Value = "  -2.21%  "

$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("D17").Value = "2.583.55"
$ws.Range("E17").Value = "  -3.75%  "

$ws.Range("E18").Value = "  -3.47%  "

$ws.Range("E19").Value = "  -0.53%  "

Set-TextValue "D20" "338.49"
$ws.Range("E20").Value = "  -1.83%  "

Set-TextValue "D22" "0.998"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("E23").Value = "  -3.09%  "

Set-TextValue "D24" "63.56"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("E25").Value = "  -0.89%  "

Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("E28").Value = "  +3.09%  "

$ws.Range("E29").Value = "  -2.54%  "

$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("E31").Value = "  -2.69%  "

Set-TextValue "D32" "162.86"
$ws.Range("E32").Value = "  -1.84%  "

Set-TextValue "D33" "4.88"
$ws.Range("E33").Value = "  +1.26%  "

$ws.Range("E34").Value = "  +0.01%  "

Set-TextValue "D35" "19.16"
$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("E36").Value = "  -1.73%  "

Set-TextValue "D37" "1.79"
$ws.Range("E37").Value = "  +0.04%  "

Set-TextValue "D38" "6.04"
$ws.Range("E38").Value = "  -1.07%  "

Set-TextValue "D39" "325.60"
$ws.Range("E39").Value = "  -4.78%  "

Set-TextValue "D40" "0.901"
$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("E41").Value = "  +0.00%  "

Set-TextValue "D42" "37.50"
$ws.Range("E42").Value = "  -1.43%  "

Set-TextValue "D43" "20.55"
$ws.Range("E43").Value = "  -0.97%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("E45").Value = "  -2.19%  "

$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("E47").Value = "  -3.05%  "

$ws.Range("E48").Value = "  -0.64%  "

Set-TextValue "D49" "19.46"
$ws.Range("E49").Value = "  -3.92%  "

$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").Value = "2.045.36"
$ws.Range("E51").Value = "  -2.12%  "
